# India Super League workbook update - 07-04-2024 22:30
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header-style (col A) and date-style (col E) formatting
# from row 122 onto the newly created rows 123-127 so no new cell styles are
# introduced into the workbook (matches the target which leaves styles.xml untouched).
$ws.Cells.Item(122, 1).Copy() | Out-Null
$ws.Range("A123:A127").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(122, 5).Copy() | Out-Null
$ws.Range("E123:E127").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 122
$ws.Cells.Item(122, 1).Value = 120
$ws.Cells.Item(122, 2).Value = 7749762
$ws.Cells.Item(122, 3).Value = "India Super League"
$ws.Cells.Item(122, 4).Value = "India Super League"
$ws.Cells.Item(122, 5).Value = 45388.35416666666
$ws.Cells.Item(122, 6).Value = "Punjab FC"
$ws.Cells.Item(122, 7).Value = "Mohun Bagan SG"
$ws.Cells.Item(122, 8).Value = 0
$ws.Cells.Item(122, 9).Value = 1
$ws.Cells.Item(122, 10).Value = "A"
$ws.Cells.Item(122, 11).Value = 4.333
$ws.Cells.Item(122, 12).Value = 3.6
$ws.Cells.Item(122, 13).Value = 1.75
$ws.Cells.Item(122, 14).Value = 4.333
$ws.Cells.Item(122, 15).Value = 4
$ws.Cells.Item(122, 16).Value = 1.727
$ws.Cells.Item(122, 17).Value = 0.75
$ws.Cells.Item(122, 18).Value = 1.875
$ws.Cells.Item(122, 19).Value = 1.975
$ws.Cells.Item(122, 20).Value = 2.75
$ws.Cells.Item(122, 21).Value = 1.825
$ws.Cells.Item(122, 22).Value = 2.025
$ws.Cells.Item(122, 23).Value = -1
$ws.Cells.Item(122, 24).Value = -1
$ws.Cells.Item(122, 25).Value = 0.7270000000000001
$ws.Cells.Item(122, 26).Value = -0.5
$ws.Cells.Item(122, 27).Value = 0.4875
$ws.Cells.Item(122, 28).Value = -1
$ws.Cells.Item(122, 29).Value = 1.025

# Row 123
$ws.Cells.Item(123, 1).Value = 121
$ws.Cells.Item(123, 2).Value = 7749471
$ws.Cells.Item(123, 3).Value = "India Super League"
$ws.Cells.Item(123, 4).Value = "India Super League"
$ws.Cells.Item(123, 5).Value = 45388.45833333334
$ws.Cells.Item(123, 6).Value = "Northeast United"
$ws.Cells.Item(123, 7).Value = "Kerala Blasters"
$ws.Cells.Item(123, 8).Value = 2
$ws.Cells.Item(123, 9).Value = 0
$ws.Cells.Item(123, 10).Value = "H"
$ws.Cells.Item(123, 11).Value = 2
$ws.Cells.Item(123, 12).Value = 3.4
$ws.Cells.Item(123, 13).Value = 3.1
$ws.Cells.Item(123, 14).Value = 1.45
$ws.Cells.Item(123, 15).Value = 4.2
$ws.Cells.Item(123, 16).Value = 5.25
$ws.Cells.Item(123, 17).Value = -1
$ws.Cells.Item(123, 18).Value = 1.8
$ws.Cells.Item(123, 19).Value = 2
$ws.Cells.Item(123, 20).Value = 3
$ws.Cells.Item(123, 21).Value = 1.925
$ws.Cells.Item(123, 22).Value = 1.875
$ws.Cells.Item(123, 23).Value = 0.45
$ws.Cells.Item(123, 24).Value = -1
$ws.Cells.Item(123, 25).Value = -1
$ws.Cells.Item(123, 26).Value = 0.8
$ws.Cells.Item(123, 27).Value = -1
$ws.Cells.Item(123, 28).Value = -1
$ws.Cells.Item(123, 29).Value = 0.875

# Row 124
$ws.Cells.Item(124, 1).Value = 122
$ws.Cells.Item(124, 2).Value = 7751762
$ws.Cells.Item(124, 3).Value = "India Super League"
$ws.Cells.Item(124, 4).Value = "India Super League"
$ws.Cells.Item(124, 5).Value = 45389.45833333334
$ws.Cells.Item(124, 6).Value = "East Bengal Club"
$ws.Cells.Item(124, 7).Value = "Bengaluru"
$ws.Cells.Item(124, 8).Value = 2
$ws.Cells.Item(124, 9).Value = 1
$ws.Cells.Item(124, 10).Value = "H"
$ws.Cells.Item(124, 11).Value = 2.55
$ws.Cells.Item(124, 12).Value = 3.2
$ws.Cells.Item(124, 13).Value = 2.55
$ws.Cells.Item(124, 14).Value = 2.45
$ws.Cells.Item(124, 15).Value = 3.2
$ws.Cells.Item(124, 16).Value = 2.7
$ws.Cells.Item(124, 17).Value = 0
$ws.Cells.Item(124, 18).Value = 1.8
$ws.Cells.Item(124, 19).Value = 2.05
$ws.Cells.Item(124, 20).Value = 2.5
$ws.Cells.Item(124, 21).Value = 1.85
$ws.Cells.Item(124, 22).Value = 2
$ws.Cells.Item(124, 23).Value = 1.45
$ws.Cells.Item(124, 24).Value = -1
$ws.Cells.Item(124, 25).Value = -1
$ws.Cells.Item(124, 26).Value = 0.8
$ws.Cells.Item(124, 27).Value = -1
$ws.Cells.Item(124, 28).Value = 0.8500000000000001
$ws.Cells.Item(124, 29).Value = -1

# Row 125
$ws.Cells.Item(125, 1).Value = 123
$ws.Cells.Item(125, 2).Value = 7749774
$ws.Cells.Item(125, 3).Value = "India Super League"
$ws.Cells.Item(125, 4).Value = "India Super League"
$ws.Cells.Item(125, 5).Value = 45390.45833333334
$ws.Cells.Item(125, 6).Value = "Mumbai City FC"
$ws.Cells.Item(125, 7).Value = "Odisha FC"
$ws.Cells.Item(125, 11).Value = 1.95
$ws.Cells.Item(125, 12).Value = 3.5
$ws.Cells.Item(125, 13).Value = 3.4
$ws.Cells.Item(125, 14).Value = 1.833
$ws.Cells.Item(125, 15).Value = 3.6
$ws.Cells.Item(125, 16).Value = 3.75
$ws.Cells.Item(125, 17).Value = -0.5
$ws.Cells.Item(125, 18).Value = 1.9
$ws.Cells.Item(125, 19).Value = 1.9
$ws.Cells.Item(125, 20).Value = 2.75
$ws.Cells.Item(125, 21).Value = 1.9
$ws.Cells.Item(125, 22).Value = 1.9
$ws.Cells.Item(125, 23).Value = 0
$ws.Cells.Item(125, 24).Value = 0
$ws.Cells.Item(125, 25).Value = 0
$ws.Cells.Item(125, 26).Value = 0
$ws.Cells.Item(125, 27).Value = 0

# Row 126
$ws.Cells.Item(126, 1).Value = 124
$ws.Cells.Item(126, 2).Value = 7751763
$ws.Cells.Item(126, 3).Value = "India Super League"
$ws.Cells.Item(126, 4).Value = "India Super League"
$ws.Cells.Item(126, 5).Value = 45391.35416666666
$ws.Cells.Item(126, 6).Value = "Jamshedpur FC"
$ws.Cells.Item(126, 7).Value = "FC Goa"
$ws.Cells.Item(126, 11).Value = 4
$ws.Cells.Item(126, 12).Value = 3.6
$ws.Cells.Item(126, 13).Value = 1.727
$ws.Cells.Item(126, 14).Value = 4
$ws.Cells.Item(126, 15).Value = 3.6
$ws.Cells.Item(126, 16).Value = 1.7
$ws.Cells.Item(126, 17).Value = 0.75
$ws.Cells.Item(126, 18).Value = 1.85
$ws.Cells.Item(126, 19).Value = 1.95
$ws.Cells.Item(126, 20).Value = 3
$ws.Cells.Item(126, 21).Value = 2
$ws.Cells.Item(126, 22).Value = 1.8
$ws.Cells.Item(126, 23).Value = 0
$ws.Cells.Item(126, 24).Value = 0
$ws.Cells.Item(126, 25).Value = 0
$ws.Cells.Item(126, 26).Value = 0
$ws.Cells.Item(126, 27).Value = 0

# Row 127
$ws.Cells.Item(127, 1).Value = 125
$ws.Cells.Item(127, 2).Value = 7751764
$ws.Cells.Item(127, 3).Value = "India Super League"
$ws.Cells.Item(127, 4).Value = "India Super League"
$ws.Cells.Item(127, 5).Value = 45391.45833333334
$ws.Cells.Item(127, 6).Value = "Chennaiyin FC"
$ws.Cells.Item(127, 7).Value = "Northeast United"
$ws.Cells.Item(127, 11).Value = 2.2
$ws.Cells.Item(127, 12).Value = 3.4
$ws.Cells.Item(127, 13).Value = 3
$ws.Cells.Item(127, 14).Value = 2.2
$ws.Cells.Item(127, 15).Value = 3.4
$ws.Cells.Item(127, 16).Value = 3
$ws.Cells.Item(127, 17).Value = -0.25
$ws.Cells.Item(127, 18).Value = 1.95
$ws.Cells.Item(127, 19).Value = 1.85
$ws.Cells.Item(127, 20).Value = 2.75
$ws.Cells.Item(127, 21).Value = 1.9
$ws.Cells.Item(127, 22).Value = 1.9
$ws.Cells.Item(127, 23).Value = 0
$ws.Cells.Item(127, 24).Value = 0
$ws.Cells.Item(127, 25).Value = 0
$ws.Cells.Item(127, 26).Value = 0
$ws.Cells.Item(127, 27).Value = 0
